$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 18 (Leve Item ID 5471)
$ws.Cells.Item(18, 8).Value = 200
$ws.Cells.Item(18, 9).Value = 200
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 200
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = 84
$ws.Cells.Item(18, 14).ClearContents()

# Row 28 (Leve Item ID 27772)
$ws.Cells.Item(28, 8).Value = 534.06665
$ws.Cells.Item(28, 9).Value = 568.4167
$ws.Cells.Item(28, 10).Value = 396.66666
$ws.Cells.Item(28, 11).Value = 568.4167
$ws.Cells.Item(28, 12).Value = 396.66666
$ws.Cells.Item(28, 13).Value = -83.41669999999999
$ws.Cells.Item(28, 14).Value = -1366.66666

# Row 29 (Leve Item ID 4575)
$ws.Cells.Item(29, 8).Value = 1024.5
$ws.Cells.Item(29, 9).Value = 27.222221
$ws.Cells.Item(29, 10).Value = 10000
$ws.Cells.Item(29, 11).Value = 81.666663
$ws.Cells.Item(29, 12).Value = 30000
$ws.Cells.Item(29, 13).Value = 199.333337
$ws.Cells.Item(29, 14).Value = -30562

# Row 31 (Leve Item ID 4576)
$ws.Cells.Item(31, 8).Value = 11156
$ws.Cells.Item(31, 9).Value = 11156
$ws.Cells.Item(31, 11).Value = 33468
$ws.Cells.Item(31, 13).Value = -33238

# Row 38 (Leve Item ID 4599)
$ws.Cells.Item(38, 8).Value = 59.52941
$ws.Cells.Item(38, 9).Value = 59.52941
$ws.Cells.Item(38, 10).Value = 0
$ws.Cells.Item(38, 11).Value = 178.58823
$ws.Cells.Item(38, 12).Value = 0
$ws.Cells.Item(38, 13).Value = 193.41177
$ws.Cells.Item(38, 14).ClearContents()

# Row 58 (Leve Item ID 4606)
$ws.Cells.Item(58, 8).Value = 1513.625
$ws.Cells.Item(58, 9).Value = 172
$ws.Cells.Item(58, 10).Value = 2557.111
$ws.Cells.Item(58, 11).Value = 516
$ws.Cells.Item(58, 12).Value = 7671.333
$ws.Cells.Item(58, 13).Value = -366
$ws.Cells.Item(58, 14).Value = -7971.333

# Row 137 (Leve Item ID 44013)
$ws.Cells.Item(137, 8).Value = 2976.6875
$ws.Cells.Item(137, 9).Value = 1498.2222
$ws.Cells.Item(137, 10).Value = 4877.5713
$ws.Cells.Item(137, 11).Value = 4494.6666
$ws.Cells.Item(137, 12).Value = 14632.7139
$ws.Cells.Item(137, 13).Value = -1944.6666
$ws.Cells.Item(137, 14).Value = -19732.7139

$ws = $wb.Worksheets.Item("ARM")
# Row 3 (Leve Item ID 2494)
$ws.Cells.Item(3, 8).Value = 22666.666
$ws.Cells.Item(3, 9).Value = 22666.666
$ws.Cells.Item(3, 11).Value = 22666.666
$ws.Cells.Item(3, 13).Value = -22551.666

# Row 5 (Leve Item ID 5091)
$ws.Cells.Item(5, 8).Value = 1000
$ws.Cells.Item(5, 9).Value = 1000
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 1000
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = -888
$ws.Cells.Item(5, 14).ClearContents()

# Row 22 (Leve Item ID 2497)
$ws.Cells.Item(22, 8).Value = 1000
$ws.Cells.Item(22, 9).Value = 1000
$ws.Cells.Item(22, 11).Value = 1000
$ws.Cells.Item(22, 13).Value = -701

# Row 28 (Leve Item ID 19534)
$ws.Cells.Item(28, 8).Value = 66185.5
$ws.Cells.Item(28, 9).Value = 66185.5
$ws.Cells.Item(28, 11).Value = 66185.5
$ws.Cells.Item(28, 13).Value = -65993.5

# Row 41 (Leve Item ID 2501)
$ws.Cells.Item(41, 8).Value = 2312.375
$ws.Cells.Item(41, 9).Value = 2312.375
$ws.Cells.Item(41, 11).Value = 2312.375
$ws.Cells.Item(41, 13).Value = -1898.375

# Row 99 (Leve Item ID 19534)
$ws.Cells.Item(99, 8).Value = 66185.5
$ws.Cells.Item(99, 9).Value = 66185.5
$ws.Cells.Item(99, 11).Value = 66185.5
$ws.Cells.Item(99, 13).Value = -63190.5

# Row 139 (Leve Item ID 42321)
$ws.Cells.Item(139, 8).Value = 58389
$ws.Cells.Item(139, 10).Value = 57927.9
$ws.Cells.Item(139, 12).Value = 57927.9
$ws.Cells.Item(139, 14).Value = -68207.89999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 4 (Leve Item ID 5091)
$ws.Cells.Item(4, 8).Value = 1000
$ws.Cells.Item(4, 9).Value = 1000
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 1000
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -885
$ws.Cells.Item(4, 14).ClearContents()

# Row 8 (Leve Item ID 2507)
$ws.Cells.Item(8, 8).Value = 12580
$ws.Cells.Item(8, 9).Value = 8444.444
$ws.Cells.Item(8, 11).Value = 8444.444
$ws.Cells.Item(8, 13).Value = -8304.444

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Cells.Item(31, 8).Value = 3925.5542
$ws.Cells.Item(31, 9).Value = 1393.1154
$ws.Cells.Item(31, 11).Value = 1393.1154
$ws.Cells.Item(31, 13).Value = -1098.1154

# Row 34 (Leve Item ID 44023)
$ws.Cells.Item(34, 8).Value = 3925.5542
$ws.Cells.Item(34, 9).Value = 1393.1154
$ws.Cells.Item(34, 11).Value = 1393.1154
$ws.Cells.Item(34, 13).Value = -1191.1154

$ws = $wb.Worksheets.Item("CUL")
# Row 3 (Leve Item ID 44094)
$ws.Cells.Item(3, 8).Value = 5590.231
$ws.Cells.Item(3, 9).Value = 2823.3333
$ws.Cells.Item(3, 10).Value = 7961.857
$ws.Cells.Item(3, 11).Value = 8469.999899999999
$ws.Cells.Item(3, 12).Value = 23885.571
$ws.Cells.Item(3, 13).Value = -8357.999899999999
$ws.Cells.Item(3, 14).Value = -24109.571

# Row 39 (Leve Item ID 4712)
$ws.Cells.Item(39, 8).Value = 5152.222
$ws.Cells.Item(39, 10).Value = 5152.222
$ws.Cells.Item(39, 12).Value = 15456.666
$ws.Cells.Item(39, 14).Value = -16044.666

# Row 64 (Leve Item ID 12861)
$ws.Cells.Item(64, 8).Value = 3299.2666
$ws.Cells.Item(64, 9).Value = 1744.8
$ws.Cells.Item(64, 10).Value = 4076.5
$ws.Cells.Item(64, 11).Value = 5234.4
$ws.Cells.Item(64, 12).Value = 12229.5
$ws.Cells.Item(64, 13).Value = -4964.4
$ws.Cells.Item(64, 14).Value = -12769.5

# Row 67 (Leve Item ID 12861)
$ws.Cells.Item(67, 8).Value = 3299.2666
$ws.Cells.Item(67, 9).Value = 1744.8
$ws.Cells.Item(67, 10).Value = 4076.5
$ws.Cells.Item(67, 11).Value = 5234.4
$ws.Cells.Item(67, 12).Value = 12229.5
$ws.Cells.Item(67, 13).Value = -4298.4
$ws.Cells.Item(67, 14).Value = -14101.5

# Row 108 (Leve Item ID 27853)
$ws.Cells.Item(108, 8).Value = 2880.6
$ws.Cells.Item(108, 9).Value = 750.2857
$ws.Cells.Item(108, 11).Value = 2250.8571
$ws.Cells.Item(108, 13).Value = 629.1428999999998

# Row 113 (Leve Item ID 27843)
$ws.Cells.Item(113, 8).Value = 244516.83
$ws.Cells.Item(113, 9).Value = 435353.47
$ws.Cells.Item(113, 10).Value = 670
$ws.Cells.Item(113, 11).Value = 1306060.41
$ws.Cells.Item(113, 12).Value = 2010
$ws.Cells.Item(113, 13).Value = -1303890.41
$ws.Cells.Item(113, 14).Value = -6350

# Row 118 (Leve Item ID 27872)
$ws.Cells.Item(118, 8).Value = 3368.6667
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 10).Value = 3368.6667
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 12).Value = 10106.0001
$ws.Cells.Item(118, 13).ClearContents()
$ws.Cells.Item(118, 14).Value = -12592.0001

# Row 122 (Leve Item ID 36078)
$ws.Cells.Item(122, 8).Value = 811.5454999999999
$ws.Cells.Item(122, 10).Value = 826.75
$ws.Cells.Item(122, 12).Value = 7440.75
$ws.Cells.Item(122, 14).Value = -12340.75

# Row 126 (Leve Item ID 36045)
$ws.Cells.Item(126, 8).Value = 2999.875
$ws.Cells.Item(126, 9).Value = 3000
$ws.Cells.Item(126, 10).Value = 2999.8572
$ws.Cells.Item(126, 11).Value = 9000
$ws.Cells.Item(126, 12).Value = 8999.571599999999
$ws.Cells.Item(126, 13).Value = -4060
$ws.Cells.Item(126, 14).Value = -18879.5716

# Row 132 (Leve Item ID 43972)
$ws.Cells.Item(132, 8).Value = 1698.44
$ws.Cells.Item(132, 9).Value = 1132.5834
$ws.Cells.Item(132, 10).Value = 2220.7693
$ws.Cells.Item(132, 11).Value = 10193.2506
$ws.Cells.Item(132, 12).Value = 19986.9237
$ws.Cells.Item(132, 13).Value = -7663.250599999999
$ws.Cells.Item(132, 14).Value = -25046.9237

# Row 133 (Leve Item ID 44073)
$ws.Cells.Item(133, 8).Value = 4320
$ws.Cells.Item(133, 9).Value = 1985.25
$ws.Cells.Item(133, 10).Value = 5756.769
$ws.Cells.Item(133, 11).Value = 5955.75
$ws.Cells.Item(133, 12).Value = 17270.307
$ws.Cells.Item(133, 13).Value = -895.75
$ws.Cells.Item(133, 14).Value = -27390.307

$ws = $wb.Worksheets.Item("LTW")
# Row 82 (Leve Item ID 12565)
$ws.Cells.Item(82, 8).Value = 1895.6666
$ws.Cells.Item(82, 9).Value = 1830.5555
$ws.Cells.Item(82, 10).Value = 1993.3334
$ws.Cells.Item(82, 11).Value = 1830.5555
$ws.Cells.Item(82, 12).Value = 1993.3334
$ws.Cells.Item(82, 13).Value = -1469.5555
$ws.Cells.Item(82, 14).Value = -2715.3334

# Row 85 (Leve Item ID 12565)
$ws.Cells.Item(85, 8).Value = 1895.6666
$ws.Cells.Item(85, 9).Value = 1830.5555
$ws.Cells.Item(85, 10).Value = 1993.3334
$ws.Cells.Item(85, 11).Value = 1830.5555
$ws.Cells.Item(85, 12).Value = 1993.3334
$ws.Cells.Item(85, 13).Value = -582.5554999999999
$ws.Cells.Item(85, 14).Value = -4489.3334

# Row 132 (Leve Item ID 44058)
$ws.Cells.Item(132, 8).Value = 4222.222
$ws.Cells.Item(132, 9).Value = 4629.2354
$ws.Cells.Item(132, 11).Value = 13887.7062
$ws.Cells.Item(132, 13).Value = -11357.7062

$ws = $wb.Worksheets.Item("WVR")
# Row 16 (Leve Item ID 26304)
$ws.Cells.Item(16, 8).Value = 45165.668
$ws.Cells.Item(16, 10).Value = 45165.668
$ws.Cells.Item(16, 12).Value = 45165.668
$ws.Cells.Item(16, 14).Value = -45749.668

